$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 42
$ws.Range("E5").Value = 138
$ws.Range("E10").Value = 592
$ws.Range("F10").Value = 300
$ws.Range("H10").Value = 396
$ws.Range("E11").Value = 376
$ws.Range("F11").Value = 204
$ws.Range("H11").Value = 268
$ws.Range("E12").Value = 583
$ws.Range("F12").Value = 323
$ws.Range("H12").Value = 409
$ws.Range("E13").Value = 143
$ws.Range("E14").Value = 128
$ws.Range("E15").Value = 175
$ws.Range("F15").Value = 75
$ws.Range("H15").Value = 125
$ws.Range("E18").Value = 53
$ws.Range("E20").Value = 90
$ws.Range("E21").Value = 141
$ws.Range("E23").Value = 206
$ws.Range("F23").Value = 101
$ws.Range("H23").Value = 152
$ws.Range("E24").Value = 222
$ws.Range("E25").Value = 284
$ws.Range("E26").Value = 161
$ws.Range("F26").Value = 99
$ws.Range("H26").Value = 124
$ws.Range("E27").Value = 340
$ws.Range("F27").Value = 181
$ws.Range("H27").Value = 262
$ws.Range("E28").Value = 204
$ws.Range("F28").Value = 84
$ws.Range("H28").Value = 136
$ws.Range("E29").Value = 173
$ws.Range("F29").Value = 102
$ws.Range("H29").Value = 143
$ws.Range("E30").Value = 219
$ws.Range("E32").Value = 189
$ws.Range("F32").Value = 116
$ws.Range("H32").Value = 154
$ws.Range("E33").Value = 305
$ws.Range("E35").Value = 156
$ws.Range("E36").Value = 77
$ws.Range("E37").Value = 168
$ws.Range("E38").Value = 95
$ws.Range("E39").Value = 182
$ws.Range("E40").Value = 270
$ws.Range("E41").Value = 404
$ws.Range("F41").Value = 196
$ws.Range("H41").Value = 288
$ws.Range("E42").Value = 393
$ws.Range("F42").Value = 219
$ws.Range("H42").Value = 280
$ws.Range("E43").Value = 125
$ws.Range("E44").Value = 320
$ws.Range("F44").Value = 164
$ws.Range("H44").Value = 232
$ws.Range("E45").Value = 153
$ws.Range("F45").Value = 75
$ws.Range("H45").Value = 114
$ws.Range("E46").Value = 336
$ws.Range("F46").Value = 189
$ws.Range("H46").Value = 252
$ws.Range("E47").Value = 472
$ws.Range("F47").Value = 249
$ws.Range("H47").Value = 341
$ws.Range("E48").Value = 225
$ws.Range("E49").Value = 294
$ws.Range("F50").Value = 122
$ws.Range("H50").Value = 193
$ws.Range("E51").Value = 247
$ws.Range("F51").Value = 111
$ws.Range("H51").Value = 185
